$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the existing data range (rows 3-173, columns A-D) ascending by column A,
# matching the re-sorted order seen in the target workbook.
$ws.Range("A3:D173").Sort($ws.Range("A3:A173"), 1)

# Append the new row (674 / Note/Donate / Date has been set / Response) at the
# very end of the table (row 174), unsorted, just like the previous batch of
# newly-added rows (now sorted in above) had been appended at the bottom.
$ws.Range("A174").Value2 = 674
$ws.Range("B174").Value2 = "Note/Donate"
$ws.Range("C174").Value2 = "Date has been set"
$ws.Range("D174").Value2 = "Response"

# Move the active selection/view down to the newly-added row, as in the
# target workbook.
$ws.Range("C174").Select()
